# Update "想去人数" (column F) figures on the 展览 and 全部类型 sheets,
# matching the data refresh recorded in the xml diff.

$wb = $excel.ActiveWorkbook

# row -> new F value, for sheet "展览"
$exhibition = @{
    2  = 72
    3  = 21713
    4  = 823
    5  = 340
    6  = 1139
    7  = 33
    8  = 8097
    10 = 51
    12 = 328
    15 = 192
    18 = 21
    23 = 58
    24 = 99
    25 = 94
    26 = 366
    27 = 1220
    31 = 614
    33 = 162
    34 = 5186
    37 = 67
    39 = 13376
    40 = 1378
    41 = 153
    42 = 62
    45 = 469
    46 = 4083
    47 = 43
}

# row -> new F value, for sheet "全部类型"
$allTypes = @{
    2  = 72
    3  = 21713
    4  = 823
    5  = 1139
    6  = 33
    7  = 8097
    9  = 51
    10 = 328
    13 = 192
    16 = 21
    20 = 58
    21 = 99
    22 = 94
    23 = 366
    24 = 1220
    29 = 614
    32 = 162
    34 = 5186
    37 = 67
    39 = 13376
    40 = 1378
    41 = 153
    42 = 62
    45 = 469
    46 = 4083
    47 = 43
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $exhibition.Keys) {
    $ws1.Range("F$row").Value = $exhibition[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypes.Keys) {
    $ws4.Range("F$row").Value = $allTypes[$row]
}
